# Automatic update of files.
# Applies the row-content permutation for rows 16-21 on the "Artfynd" sheet:
#   new row 16 <= old row 20
#   new row 17 <= old row 16
#   new row 18 <= old row 17
#   new row 19 <= old row 21
#   new row 20 <= old row 18
#   new row 21 <= old row 19
# Implemented as direct per-cell updates (only the cells whose value actually
# changes are touched), matching the published diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

$ws.Range("A16").Value = 91841297
$ws.Range("B16").Value = 56395
$ws.Range("E16").Value = 100109
$ws.Range("F16").Value = 'Tretåig hackspett'
$ws.Range("G16").Value = 'Picoides tridactylus'
$ws.Range("Q16").Value = 435165.9665010545
$ws.Range("R16").Value = 6971255.90823228
$ws.Range("AC16").Value = 'Ringhack'
$ws.Range("A17").Value = 91841152
$ws.Range("B17").Value = 56411
$ws.Range("D17").Value = 'NT'
$ws.Range("E17").Value = 100049
$ws.Range("F17").Value = 'Spillkråka'
$ws.Range("G17").Value = 'Dryocopus martius'
$ws.Range("H17").Value = '(Linnaeus, 1758)'
$ws.Range("Q17").Value = 435739.0389247189
$ws.Range("R17").Value = 6971194.847051037
$ws.Range("AC17").Value = 'Födosökshack'
$ws.Range("A18").Value = 91841302
$ws.Range("B18").Value = 96354
$ws.Range("D18").Value = 'LC'
$ws.Range("E18").Value = 221952
$ws.Range("F18").Value = 'Spindelblomster'
$ws.Range("G18").Value = 'Neottia cordata'
$ws.Range("H18").Value = '(L.) Rich.'
$ws.Range("K18").ClearContents()
$ws.Range("L18").ClearContents()
$ws.Range("M18").ClearContents()
$ws.Range("P18").Value = 'Jämtlands län, Jmt'
$ws.Range("Q18").Value = 435430.9732251486
$ws.Range("R18").Value = 6971269.849731034
$ws.Range("S18").Value = 25
$ws.Range("Y18").Value = '''2019-06-01'
$ws.Range("AA18").Value = '''2019-10-31'
$ws.Range("AW18").Value = 'Erland Lindblad'
$ws.Range("AX18").Value = 'Via Erland Lindblad'
$ws.Range("AY18").Value = 'Kontinuitetsskogar och skogar med höga naturvärden ovan och i nära anslutning till fjällnära gränsen'
$ws.Range("A19").Value = 106029917
$ws.Range("Q19").Value = 435518.6180496986
$ws.Range("R19").Value = 6970728.401884655
$ws.Range("A20").Value = 106029936
$ws.Range("K20").Value = 'adult'
$ws.Range("L20").Value = 'honfärgad'
$ws.Range("M20").Value = 'födosökande'
$ws.Range("P20").Value = 'Stavbrännhöjden, Jmt'
$ws.Range("Q20").Value = 435672.668452872
$ws.Range("R20").Value = 6971354.125670544
$ws.Range("S20").Value = 10
$ws.Range("Y20").Value = '''2023-01-16'
$ws.Range("AA20").Value = '''2023-01-16'
$ws.Range("AC20").ClearContents()
$ws.Range("AW20").Value = 'Benny Öwre'
$ws.Range("AX20").Value = 'Benny Öwre'
$ws.Range("AY20").ClearContents()
$ws.Range("A21").Value = 106029937
$ws.Range("Q21").Value = 435277.8438721292
$ws.Range("R21").Value = 6971109.868824044
